$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. DoPlaySheet -> RVL.DoPlayTest (rename Action + Param Name for every RVL play-script row)
#    and retarget the Param Value path from the *.rvl.xlsx script to the matching *.sstest file
#    2. TestPrepare with navigation to site according to Config.xlsx (new DataSources/DataOrigin/DataOutput rows)

# Row 3 - DefaultLogins
$ws.Range("D3").Value = "DoPlayTest"
$ws.Range("E3").Value = "sstestPath"
$ws.Range("G3").Value = "%WORKDIR%\DefaultLogins\DefaultLogins.sstest"

# Row 4 - CheckDashboardTiles
$ws.Range("D4").Value = "DoPlayTest"
$ws.Range("E4").Value = "sstestPath"
$ws.Range("G4").Value = "%WORKDIR%\CheckDashboardTiles\CheckDashboardTiles.sstest"

# Row 5 - RegisterPatient
$ws.Range("D5").Value = "DoPlayTest"
$ws.Range("E5").Value = "sstestPath"
$ws.Range("G5").Value = "%WORKDIR%\RegisterPatient\RegisterPatient.sstest"

# Row 7 - DataSources
$ws.Range("D7").Value = "DoPlayTest"
$ws.Range("E7").Value = "sstestPath"
$ws.Range("G7").Value = "%WORKDIR%\DataSources\DataSources.sstest"

# Row 8 - DataOrigin
$ws.Range("D8").Value = "DoPlayTest"
$ws.Range("E8").Value = "sstestPath"
$ws.Range("G8").Value = "%WORKDIR%\DataOrigin\DataOrigin.sstest"

# Row 9 - DataOutput
$ws.Range("D9").Value = "DoPlayTest"
$ws.Range("E9").Value = "sstestPath"
$ws.Range("G9").Value = "%WORKDIR%\DataOutput\DataOutput.sstest"

# Remove the blank separator row 6 so DataSources/DataOrigin/DataOutput move directly
# below the RegisterPatient block (rows shift up by one, dimension becomes A1:H33)
$ws.Rows("6:6").Delete()
